$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Occupied" header and values in column D
$ws.Range("D1").Value = "Occupied"
$ws.Range("D2").Value = "y"
$ws.Range("D3").Value = "y"
$ws.Range("D4").Value = "n"
$ws.Range("D5").Value = "n"
$ws.Range("D6").Value = "y"
$ws.Range("D7").Value = "n"
$ws.Range("D8").Value = "y"
$ws.Range("D9").Value = "y"

# Match the new row height for the header row
$ws.Rows.Item(1).RowHeight = 12.5

# Update the selected cell
$ws.Range("H12").Select() | Out-Null
